$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header changes: insert AWC at E1, shift TVN to F1, add CTC at G1 ---
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("E1").Value = "AWC"
$ws.Range("F1").Value = "TVN"
$ws.Range("G1").Value = "CTC"

# --- Data arrays (rows 2-51) ---
$Avals = @(7, 189, 309, 335, 381, 384, 388, 391, 395, 400, 403, 407, 425, 428, 432, 435, 440, 443, 448, 451, 454, 461, 464, 470, 473, 477, 481, 488, 496, 499, 502, 508, 512, 518, 568, 576, 581, 983, 999, 1043, 1148, 1383, 1387, 1397, 1458, 1547, 1551, 1558, 1561, 1608)
$Bvals = @("5:28 AM", "7:00 AM", "7:59 AM", "8:12 AM", "8:35 AM", "8:37 AM", "8:39 AM", "8:40 AM", "8:42 AM", "8:45 AM", "8:47 AM", "8:48 AM", "8:57 AM", "8:59 AM", "9:01 AM", "9:03 AM", "9:05 AM", "9:06 AM", "9:09 AM", "9:11 AM", "9:12 AM", "9:15 AM", "9:17 AM", "9:20 AM", "9:21 AM", "9:24 AM", "9:25 AM", "9:29 AM", "9:33 AM", "9:34 AM", "9:36 AM", "9:39 AM", "9:41 AM", "9:44 AM", "10:09 A", "10:13 A", "10:15 A", "1:36 PM", "1:44 PM", "2:06 PM", "2:59 PM", "4:56 PM", "4:58 PM", "5:03 PM", "5:34 PM", "6:19 PM", "6:20 PM", "6:24 PM", "6:25 PM", "6:49 PM")
$Cvals = @(202.2, 5668.34, 9248.33, 10035.655, 11413.92, 11516.145, 11611.81, 11704.12, 11822.4, 11997.16, 12088.72, 12193.74, 12725.265, 12831.085, 12941.465, 13048.17, 13177.29, 13276.226667, 13418.07, 13527.81, 13619.07, 13803.455, 13906.515, 14073.82, 14174.17, 14308.02, 14418.09, 14621.995, 14854.345, 14960.895, 15056.18, 15216.34, 15345.09, 15517.59, 17019.785, 17250.1, 17413.61, 29464.66, 29954.343333, 31284.45, 34414.095, 41472.38, 41596.02, 41887.765, 43728.35, 46409.81, 46502.28, 46717.07, 46813.12, 48231.23)
$Dvals = @(232.2, 5698.34, 9278.33, 10065.655, 11443.92, 11546.145, 11641.81, 11734.12, 11852.4, 12027.16, 12118.72, 12223.74, 12755.265, 12861.085, 12971.465, 13078.17, 13207.29, 13306.226667, 13448.07, 13557.81, 13649.07, 13833.455, 13936.515, 14103.82, 14204.17, 14338.02, 14448.09, 14651.995, 14884.345, 14990.895, 15086.18, 15246.34, 15375.09, 15547.59, 17049.785, 17280.1, 17443.61, 29494.66, 29984.343333, 31314.45, 34444.095, 41502.38, 41626.02, 41917.765, 43758.35, 46439.81, 46532.28, 46747.07, 46843.12, 48261.23)
$Evals = @(11.26, 0, 0, 4.215, 1.94, 3.955, 0, 0, 0, 0, 0, 1.355, 4.76, 0, 2.38, 0, 0, 1.58, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1.415, 0, 0, 6.085, 0.775, 3.238333, 20.516667, 8.845000000000001, 0, 0, 0.803333, 3.856667, 0, 0, 7.295, 0, 2.37, 3.405, 0, 39.27, 0, 0, 0)
$Fvals = @(1.29, 1, 1.065, 1.55, 19.395, 28.035, 59.97, 20.83, 8.73, 2.8, 32.41, 58.128333, 43.895, 51.785, 48.73, 7.34, 60.2, 5.873333, 17.015, 40.62, 49.685, 45.39, 11.275, 30.935, 4.846667, 4.635, 3.05, 113.915, 19.64, 18.025, 24.46, 20.27, 2.291667, 0.58, 1.025, 3.83, 676.14, 0.78, 0.826667, 1.12, 0.8, 1.325, 1.94, 0.67, 0.5649999999999999, 1.36, 1.01, 0.64, 1.03, 0.695)
$Gvals = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.5, 0, 0, 0.333333, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.5, 0, 0.166667, 0, 1, 0, 0, 0, 0.666667, 0, 0, 0, 0, 0.5, 0, 0, 1, 0, 0, 0)

for ($i = 0; $i -lt 50; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $Avals[$i]
    $ws.Cells.Item($r, 2).Value = $Bvals[$i]
    $ws.Cells.Item($r, 3).Value = $Cvals[$i]
    $ws.Cells.Item($r, 4).Value = $Dvals[$i]
    $ws.Cells.Item($r, 5).Value = $Evals[$i]
    $ws.Cells.Item($r, 6).Value = $Fvals[$i]
    $ws.Cells.Item($r, 7).Value = $Gvals[$i]
}
